# Insert a new "Date Created (Year)*" column before the existing column C
# (description), shifting description..AUTHOR-role columns one to the
# right, and populate the new column with a year value for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at C; existing C:R shift to D:S.
$ws.Columns.Item(3).Insert()

# Header for the new column.
$ws.Range("C1").Value = "Date Created (Year)*"

# Populate the year value for each data row (numeric, not text).
$ws.Range("C2").Value = 2000
$ws.Range("C3").Value = 2000
$ws.Range("C4").Value = 2000

# New font (explicit black RGB rather than the theme color) applied to the
# freshly populated cells.
$ws.Range("C2:C4").Font.Color = 0

# Reflect the new selection state recorded in the sheet view.
$ws.Range("C1:C4").Select() | Out-Null
